# Insert two new price-report rows (weekly update) at the top of the data
# block (rows 15-16), pushing the existing rows 15-63 down to 17-65.
# Dimension grows from A1:R63 to A1:R65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before the current row 15 - this shifts rows 15:63
# down to 17:65 and keeps all their existing content/formatting intact.
$ws.Rows("15:16").Insert()

# --- New row 15: Perejil, Primera, 2023-05-25 ---
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 45071
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 100112044
$ws.Range("G15").Value = "Perejil"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 1200
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1200
$ws.Range("N15").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1200
$ws.Range("Q15").Value = 1
$ws.Range("R15").Value = "Hortaliza"

# --- New row 16: Perejil, Segunda, 2023-05-25 ---
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 45071
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 100112044
$ws.Range("G16").Value = "Perejil"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Segunda"
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 1000
$ws.Range("N16").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 1000
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = "Hortaliza"
